$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 260, pushing the existing rows
# 260-296 down to 261-297 (same as Excel's Rows.Insert on a full row).
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new weekly price record.
$ws.Range("A260").Value = 8
$ws.Range("B260").Value = "Terminal La Palmera de La Serena"
$ws.Range("C260").Value = "Coquimbo"
$ws.Range("D260").Value = 44776
$ws.Range("D260").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E260").Value = 4
$ws.Range("F260").Value = 100112003
$ws.Range("G260").Value = "Ajo"
$ws.Range("H260").Value = "Chino"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 440
$ws.Range("K260").Value = 26000
$ws.Range("L260").Value = 27000
$ws.Range("M260").Value = 26500
$ws.Range("N260").Value = "`$/caja 10 kilos"
$ws.Range("O260").Value = "China"
$ws.Range("P260").Value = 2650
$ws.Range("Q260").Value = 10
$ws.Range("R260").Value = "Hortaliza"
